# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns (J, K, L, R) for the
# 8f98d553-1ce3-480c-89cf-0dd615b9b3e6 row (row 5) on both the "zh-cn" and
# "de-de" report sheets, widens columns J and K, and adds a hyperlink on
# the newly-populated "Latest Target File" cell (mirroring the existing
# "Source File Name" hyperlink for that row).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/b9f8a019ffddfc29c1799976b534962e71f7e335/e2e/8f98d553-1ce3-480c-89cf-0dd615b9b3e6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/f9e6529fdca62dc5ea17d3dc1b9018a431f07841/e2e/8f98d553-1ce3-480c-89cf-0dd615b9b3e6.md."
$targetFileName = "8f98d553-1ce3-480c-89cf-0dd615b9b3e6.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/f9e6529fdca62dc5ea17d3dc1b9018a431f07841/e2e/8f98d553-1ce3-480c-89cf-0dd615b9b3e6.md"

function Update-HandbackRow {
    param($ws, [string]$handbackFile, [string]$handbackDate)

    # Widen the "Latest Target File" / "Latest Handback File" columns (J/K).
    $ws.Columns.Item(10).ColumnWidth = 39.1
    $ws.Columns.Item(11).ColumnWidth = 39.1

    # K5: Latest Handback File
    $ws.Range("K5").Value = $handbackFile

    # L5: Latest Handback DateTime
    $ws.Range("L5").Value = $handbackDate

    # R5: Error Detail
    $ws.Range("R5").Value = $errorDetail

    # J5: Latest Target File - add as a hyperlink (same target as the
    # "Source File Name" hyperlink for this row), which also sets the
    # display text and applies hyperlink formatting.
    $ws.Hyperlinks.Add($ws.Range("J5"), $targetUrl, "", "", $targetFileName)
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "8f98d553-1ce3-480c-89cf-0dd615b9b3e6.101513ec6eea8c5e77a7de649f9cd2268629cbb1.zh-cn.xlf" "2017-02-28 07:33:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "8f98d553-1ce3-480c-89cf-0dd615b9b3e6.101513ec6eea8c5e77a7de649f9cd2268629cbb1.de-de.xlf" "2017-02-28 07:33:30"
